# Auto-generated edit script: update market-price derived cells per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 125.3125
$ws.Range("I12").Value = 120.4
$ws.Range("K12").Value = 120.4
$ws.Range("M12").Value = 49.59999999999999
$ws.Range("H19").Value = 1138.8
$ws.Range("I19").Value = 1247.5
$ws.Range("J19").Value = 1066.3334
$ws.Range("K19").Value = 1247.5
$ws.Range("L19").Value = 1066.3334
$ws.Range("M19").Value = -1072.5
$ws.Range("N19").Value = -1416.3334
$ws.Range("H33").Value = 551.6087
$ws.Range("I33").Value = 274.1875
$ws.Range("K33").Value = 274.1875
$ws.Range("M33").Value = -45.1875
$ws.Range("H43").Value = 4539
$ws.Range("I43").Value = 3830
$ws.Range("J43").Value = 6666
$ws.Range("K43").Value = 3830
$ws.Range("L43").Value = 6666
$ws.Range("M43").Value = -3761
$ws.Range("N43").Value = -6804
$ws.Range("H116").Value = 6033.077
$ws.Range("I116").Value = 6804.636
$ws.Range("K116").Value = 6804.636
$ws.Range("M116").Value = -3362.636
$ws.Range("H137").Value = 3337.6316
$ws.Range("I137").Value = 2650.2
$ws.Range("K137").Value = 7950.599999999999
$ws.Range("M137").Value = -5400.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 825.0577
$ws.Range("I2").Value = 787.36957
$ws.Range("J2").Value = 1114
$ws.Range("K2").Value = 787.36957
$ws.Range("L2").Value = 1114
$ws.Range("M2").Value = -674.36957
$ws.Range("N2").Value = -1340
$ws.Range("H45").Value = 1668
$ws.Range("I45").Value = 1631.5555
$ws.Range("J45").Value = 1750
$ws.Range("K45").Value = 1631.5555
$ws.Range("L45").Value = 1750
$ws.Range("M45").Value = -1254.5555
$ws.Range("N45").Value = -2504
$ws.Range("H116").Value = 825.0577
$ws.Range("I116").Value = 787.36957
$ws.Range("J116").Value = 1114
$ws.Range("K116").Value = 787.36957
$ws.Range("L116").Value = 1114
$ws.Range("M116").Value = 1506.63043
$ws.Range("N116").Value = -5702

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 825.0577
$ws.Range("I3").Value = 787.36957
$ws.Range("J3").Value = 1114
$ws.Range("K3").Value = 787.36957
$ws.Range("L3").Value = 1114
$ws.Range("M3").Value = -673.36957
$ws.Range("N3").Value = -1342

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 788.5
$ws.Range("I22").Value = 466
$ws.Range("K22").Value = 466
$ws.Range("M22").Value = -116
$ws.Range("H58").Value = 4299.625
$ws.Range("I58").Value = 5045.8335
$ws.Range("K58").Value = 5045.8335
$ws.Range("M58").Value = -4842.8335
$ws.Range("H110").Value = 89999.89
$ws.Range("J110").Value = 89999.89
$ws.Range("L110").Value = 89999.89
$ws.Range("N110").Value = -98179.89
$ws.Range("H116").Value = 62500
$ws.Range("J116").Value = 62500
$ws.Range("L116").Value = 62500
$ws.Range("N116").Value = -71678
$ws.Range("H117").Value = 54856
$ws.Range("J117").Value = 54856
$ws.Range("L117").Value = 54856
$ws.Range("N117").Value = -64034
$ws.Range("H122").Value = 2575.8147
$ws.Range("I122").Value = 2881.9524
$ws.Range("J122").Value = 1504.3334
$ws.Range("K122").Value = 8645.8572
$ws.Range("L122").Value = 4513.0002
$ws.Range("M122").Value = -6195.8572
$ws.Range("N122").Value = -9413.0002
$ws.Range("H136").Value = 4299.625
$ws.Range("I136").Value = 5045.8335
$ws.Range("K136").Value = 15137.5005
$ws.Range("M136").Value = -12587.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1974
$ws.Range("I3").Value = 1974
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5922
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -5810
$ws.Range("N3").ClearContents()
$ws.Range("H18").Value = 654.8570999999999
$ws.Range("I18").Value = 514
$ws.Range("J18").Value = 1500
$ws.Range("K18").Value = 1542
$ws.Range("L18").Value = 4500
$ws.Range("M18").Value = -1373
$ws.Range("N18").Value = -4838
$ws.Range("H107").Value = 1114.25
$ws.Range("I107").Value = 949.5
$ws.Range("K107").Value = 2848.5
$ws.Range("M107").Value = -928.5
$ws.Range("H114").Value = 19599.6
$ws.Range("I114").Value = 3998
$ws.Range("J114").Value = 23500
$ws.Range("K114").Value = 11994
$ws.Range("L114").Value = 70500
$ws.Range("M114").Value = -8740
$ws.Range("N114").Value = -77008
$ws.Range("H118").Value = 1376.3334
$ws.Range("I118").Value = 1376.3334
$ws.Range("K118").Value = 4129.0002
$ws.Range("M118").Value = -2886.0002
$ws.Range("H126").Value = 675
$ws.Range("I126").Value = 675
$ws.Range("K126").Value = 2025
$ws.Range("M126").Value = 2915
$ws.Range("H131").Value = 33427.23
$ws.Range("I131").Value = 278509
$ws.Range("J131").Value = 1803.7742
$ws.Range("K131").Value = 835527
$ws.Range("L131").Value = 5411.3226
$ws.Range("M131").Value = -830487
$ws.Range("N131").Value = -15491.3226
$ws.Range("H133").Value = 6681.8184
$ws.Range("J133").Value = 7250
$ws.Range("L133").Value = 21750
$ws.Range("N133").Value = -31870
$ws.Range("H136").Value = 6833.3335
$ws.Range("I136").Value = 6833.3335
$ws.Range("K136").Value = 20500.0005
$ws.Range("M136").Value = -15400.0005
$ws.Range("H138").Value = 1111.8
$ws.Range("I138").Value = 1111.8
$ws.Range("K138").Value = 3335.4
$ws.Range("M138").Value = 1804.6
$ws.Range("H139").Value = 5178.5625
$ws.Range("I139").Value = 3077.0715
$ws.Range("K139").Value = 9231.2145
$ws.Range("M139").Value = -4091.2145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8511.235000000001
$ws.Range("J80").Value = 11259.4
$ws.Range("L80").Value = 11259.4
$ws.Range("N80").Value = -13255.4
$ws.Range("H83").Value = 8511.235000000001
$ws.Range("J83").Value = 11259.4
$ws.Range("L83").Value = 56297
$ws.Range("N83").Value = -66281
$ws.Range("H126").Value = 4034.238
$ws.Range("J126").Value = 4527.385
$ws.Range("L126").Value = 13582.155
$ws.Range("N126").Value = -18522.155
$ws.Range("H132").Value = 3261.88
$ws.Range("J132").Value = 2725.5715
$ws.Range("L132").Value = 8176.7145
$ws.Range("N132").Value = -13236.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4898.5713
$ws.Range("I7").Value = 5459.6
$ws.Range("K7").Value = 5459.6
$ws.Range("M7").Value = -5347.6
$ws.Range("H22").Value = 9066.333000000001
$ws.Range("I22").Value = 5599.5
$ws.Range("K22").Value = 5599.5
$ws.Range("M22").Value = -5304.5
$ws.Range("H27").Value = 9066.333000000001
$ws.Range("I27").Value = 5599.5
$ws.Range("K27").Value = 5599.5
$ws.Range("M27").Value = -5492.5
$ws.Range("H61").Value = 9997.5
$ws.Range("I61").Value = 9997
$ws.Range("K61").Value = 9997
$ws.Range("M61").Value = -9795
$ws.Range("H113").Value = 9997.5
$ws.Range("I113").Value = 9997
$ws.Range("K113").Value = 9997
$ws.Range("M113").Value = -7827
$ws.Range("H126").Value = 4898.5713
$ws.Range("I126").Value = 5459.6
$ws.Range("K126").Value = 16378.8
$ws.Range("M126").Value = -13908.8
$ws.Range("H132").Value = 3678.6956
$ws.Range("I132").Value = 3475.2
$ws.Range("K132").Value = 10425.6
$ws.Range("M132").Value = -7895.599999999999
$ws.Range("H136").Value = 2736.8667
$ws.Range("J136").Value = 2363.5
$ws.Range("L136").Value = 7090.5
$ws.Range("N136").Value = -12190.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1083.4546
$ws.Range("I113").Value = 692.25
$ws.Range("J113").Value = 1307
$ws.Range("K113").Value = 2076.75
$ws.Range("L113").Value = 3921
$ws.Range("M113").Value = 93.25
$ws.Range("N113").Value = -8261
$ws.Range("H122").Value = 5251.684
$ws.Range("I122").Value = 5052.6665
$ws.Range("K122").Value = 15157.9995
$ws.Range("M122").Value = -12707.9995
